$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1544.8306
$ws.Range("I15").Value = 1544.8306
$ws.Range("K15").Value = 4634.4918
$ws.Range("M15").Value = -4465.4918
$ws.Range("H70").Value = 100006200
$ws.Range("I70").Value = 2997.5
$ws.Range("J70").Value = 166674990
$ws.Range("K70").Value = 8992.5
$ws.Range("L70").Value = 500024970
$ws.Range("M70").Value = -8722.5
$ws.Range("N70").Value = -500025510
$ws.Range("H73").Value = 100006200
$ws.Range("I73").Value = 2997.5
$ws.Range("J73").Value = 166674990
$ws.Range("K73").Value = 8992.5
$ws.Range("L73").Value = 500024970
$ws.Range("M73").Value = -8056.5
$ws.Range("N73").Value = -500026842
$ws.Range("H75").Value = 272766820
$ws.Range("J75").Value = 272766820
$ws.Range("L75").Value = 272766820
$ws.Range("N75").Value = -272768692
$ws.Range("H78").Value = 272766820
$ws.Range("J78").Value = 272766820
$ws.Range("L78").Value = 818300460
$ws.Range("N78").Value = -818309820
$ws.Range("H141").Value = 5496.1304
$ws.Range("I141").Value = 3416.4443
$ws.Range("K141").Value = 10249.3329
$ws.Range("M141").Value = -5069.332900000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1919.45
$ws.Range("I2").Value = 821.61536
$ws.Range("K2").Value = 821.61536
$ws.Range("M2").Value = -708.61536
$ws.Range("H5").Value = 84
$ws.Range("I5").Value = 84
$ws.Range("K5").Value = 84
$ws.Range("M5").Value = 28
$ws.Range("H74").Value = 3045.8975
$ws.Range("I74").Value = 2819.2693
$ws.Range("K74").Value = 2819.2693
$ws.Range("M74").Value = -1945.2693
$ws.Range("H77").Value = 3045.8975
$ws.Range("I77").Value = 2819.2693
$ws.Range("K77").Value = 14096.3465
$ws.Range("M77").Value = -9728.3465
$ws.Range("H101").Value = 31474.5
$ws.Range("J101").Value = 31474.5
$ws.Range("L101").Value = 31474.5
$ws.Range("N101").Value = -37964.5
$ws.Range("H116").Value = 1919.45
$ws.Range("I116").Value = 821.61536
$ws.Range("K116").Value = 821.61536
$ws.Range("M116").Value = 1472.38464
$ws.Range("H132").Value = 1706.1025
$ws.Range("I132").Value = 1706.1025
$ws.Range("K132").Value = 5118.3075
$ws.Range("M132").Value = -2588.3075
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1919.45
$ws.Range("I3").Value = 821.61536
$ws.Range("K3").Value = 821.61536
$ws.Range("M3").Value = -707.61536
$ws.Range("H4").Value = 84
$ws.Range("I4").Value = 84
$ws.Range("K4").Value = 84
$ws.Range("M4").Value = 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1833.65
$ws.Range("I58").Value = 1871.2572
$ws.Range("K58").Value = 1871.2572
$ws.Range("M58").Value = -1668.2572
$ws.Range("H132").Value = 2988.1538
$ws.Range("J132").Value = 7167.8335
$ws.Range("L132").Value = 21503.5005
$ws.Range("N132").Value = -26563.5005
$ws.Range("H136").Value = 1833.65
$ws.Range("I136").Value = 1871.2572
$ws.Range("K136").Value = 5613.7716
$ws.Range("M136").Value = -3063.7716
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6053.222
$ws.Range("I3").Value = 6121
$ws.Range("J3").Value = 5999
$ws.Range("K3").Value = 18363
$ws.Range("L3").Value = 17997
$ws.Range("M3").Value = -18251
$ws.Range("N3").Value = -18221
$ws.Range("H4").Value = 37286080
$ws.Range("J4").Value = 17571464
$ws.Range("L4").Value = 52714392
$ws.Range("N4").Value = -52714616
$ws.Range("H26").Value = 1458.9286
$ws.Range("I26").Value = 1357.091
$ws.Range("J26").Value = 1832.3334
$ws.Range("K26").Value = 4071.273
$ws.Range("L26").Value = 5497.0002
$ws.Range("M26").Value = -3783.273
$ws.Range("N26").Value = -6073.0002
$ws.Range("H86").Value = 914.5
$ws.Range("J86").Value = 1040.5
$ws.Range("L86").Value = 3121.5
$ws.Range("N86").Value = -5493.5
$ws.Range("H89").Value = 914.5
$ws.Range("J89").Value = 1040.5
$ws.Range("L89").Value = 9364.5
$ws.Range("N89").Value = -21220.5
$ws.Range("H103").Value = 50
$ws.Range("I103").Value = 50
$ws.Range("K103").Value = 150
$ws.Range("M103").Value = 729
$ws.Range("H122").Value = 1788.4
$ws.Range("I122").Value = 1120.6666
$ws.Range("J122").Value = 1955.3334
$ws.Range("K122").Value = 10085.9994
$ws.Range("L122").Value = 17598.0006
$ws.Range("M122").Value = -7635.999400000001
$ws.Range("N122").Value = -22498.0006
$ws.Range("H131").Value = 3336464.5
$ws.Range("J131").Value = 3925111.5
$ws.Range("L131").Value = 11775334.5
$ws.Range("N131").Value = -11785414.5
$ws.Range("H133").Value = 9898.526
$ws.Range("I133").Value = 3259
$ws.Range("J133").Value = 14727.272
$ws.Range("K133").Value = 9777
$ws.Range("L133").Value = 44181.81600000001
$ws.Range("M133").Value = -4717
$ws.Range("N133").Value = -54301.81600000001
$ws.Range("H140").Value = 2510.1667
$ws.Range("I140").Value = 2540.1765
$ws.Range("K140").Value = 7620.529500000001
$ws.Range("M140").Value = -2440.529500000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26437.285
$ws.Range("J57").Value = 29176.834
$ws.Range("L57").Value = 29176.834
$ws.Range("N57").Value = -30816.834
$ws.Range("H93").Value = 20199.2
$ws.Range("J93").Value = 20199.2
$ws.Range("L93").Value = 20199.2
$ws.Range("N93").Value = -23943.2
$ws.Range("H132").Value = 23811584
$ws.Range("I132").Value = 32259930
$ws.Range("K132").Value = 96779790
$ws.Range("M132").Value = -96777260
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2572.7368
$ws.Range("I22").Value = 1649
$ws.Range("J22").Value = 3599.111
$ws.Range("K22").Value = 1649
$ws.Range("L22").Value = 3599.111
$ws.Range("M22").Value = -1354
$ws.Range("N22").Value = -4189.111
$ws.Range("H27").Value = 2572.7368
$ws.Range("I27").Value = 1649
$ws.Range("J27").Value = 3599.111
$ws.Range("K27").Value = 1649
$ws.Range("L27").Value = 3599.111
$ws.Range("M27").Value = -1542
$ws.Range("N27").Value = -3813.111
$ws.Range("H46").Value = 890.62964
$ws.Range("I46").Value = 431.7619
$ws.Range("J46").Value = 2496.6667
$ws.Range("K46").Value = 431.7619
$ws.Range("L46").Value = 2496.6667
$ws.Range("M46").Value = -243.7619
$ws.Range("N46").Value = -2872.6667
$ws.Range("H62").Value = 400021440
$ws.Range("J62").Value = 500019230
$ws.Range("L62").Value = 500019230
$ws.Range("N62").Value = -500020478
$ws.Range("H65").Value = 400021440
$ws.Range("J65").Value = 500019230
$ws.Range("L65").Value = 1500057690
$ws.Range("N65").Value = -1500063930
$ws.Range("H82").Value = 32260454
$ws.Range("I82").Value = 812.4666999999999
$ws.Range("J82").Value = 62503868
$ws.Range("K82").Value = 812.4666999999999
$ws.Range("L82").Value = 62503868
$ws.Range("M82").Value = -451.4666999999999
$ws.Range("N82").Value = -62504590
$ws.Range("H85").Value = 32260454
$ws.Range("I85").Value = 812.4666999999999
$ws.Range("J85").Value = 62503868
$ws.Range("K85").Value = 812.4666999999999
$ws.Range("L85").Value = 62503868
$ws.Range("M85").Value = 435.5333000000001
$ws.Range("N85").Value = -62506364
$ws.Range("H132").Value = 2139849.5
$ws.Range("I132").Value = 2418312.5
$ws.Range("K132").Value = 7254937.5
$ws.Range("M132").Value = -7252407.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23544.455
$ws.Range("I54").Value = 14332.333
$ws.Range("J54").Value = 26999
$ws.Range("K54").Value = 14332.333
$ws.Range("L54").Value = 26999
$ws.Range("M54").Value = -13812.333
$ws.Range("N54").Value = -28039
$ws.Range("H75").Value = 800010000
$ws.Range("J75").Value = 800010000
$ws.Range("L75").Value = 800010000
$ws.Range("N75").Value = -800011872
$ws.Range("H78").Value = 800010000
$ws.Range("J78").Value = 800010000
$ws.Range("L78").Value = 2400030000
$ws.Range("N78").Value = -2400039360
$ws.Range("H81").Value = 2933.7144
$ws.Range("I81").Value = 2766.4
$ws.Range("K81").Value = 5532.8
$ws.Range("M81").Value = -4471.8
$ws.Range("H84").Value = 2933.7144
$ws.Range("I84").Value = 2766.4
$ws.Range("K84").Value = 27664
$ws.Range("M84").Value = -22360
$ws.Range("H107").Value = 27778532
$ws.Range("J107").Value = 62500740
$ws.Range("L107").Value = 187502220
$ws.Range("N107").Value = -187506060
$ws.Range("H124").Value = 57266.832
$ws.Range("J124").Value = 57266.832
$ws.Range("L124").Value = 57266.832
$ws.Range("N124").Value = -67086.83199999999
$ws.Range("H136").Value = 4731.722
$ws.Range("I136").Value = 2670.5715
$ws.Range("J136").Value = 6043.364
$ws.Range("K136").Value = 8011.7145
$ws.Range("L136").Value = 18130.092
$ws.Range("M136").Value = -5461.7145
$ws.Range("N136").Value = -23230.092
